$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2: updated odds values ---
$ws.Range("G2").Value = 2.35
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = 3.3
$ws.Range("J2").Value = 3.2
$ws.Range("L2").Value = 4
$ws.Range("M2").Value = 1.11
$ws.Range("N2").Value = 6.5
$ws.Range("AC2").Value = 6.5
$ws.Range("AD2").Value = 6
$ws.Range("AG2").Value = 7
$ws.Range("AJ2").Value = 34
$ws.Range("AN2").Value = 4.33
$ws.Range("AU2").Value = 9.5
$ws.Range("AZ2").Value = 67

# --- Row 4: new match row appended below existing data ---
$ws.Range("A4").Value = "M5Xw3O3j"
$ws.Range("B4").Value = "18/11/2024"
$ws.Range("C4").Value = "19:00"
$ws.Range("D4").Value = "URUGUAY - PRIMERA DIVISION"
$ws.Range("E4").Value = "CA Cerro"
$ws.Range("F4").Value = "Boston River"
$ws.Range("G4").Value = 3.25
$ws.Range("H4").Value = 3.25
$ws.Range("I4").Value = 2.25
$ws.Range("J4").Value = 4
$ws.Range("K4").Value = 2.05
$ws.Range("L4").Value = 3
$ws.Range("M4").Value = 1.07
$ws.Range("N4").Value = 9
$ws.Range("O4").Value = 1.36
$ws.Range("P4").Value = 3
$ws.Range("Q4").Value = 2.15
$ws.Range("R4").Value = 1.67
$ws.Range("S4").Value = 1.5
$ws.Range("T4").Value = 2.5
$ws.Range("U4").Value = 1.91
$ws.Range("V4").Value = 1.8
$ws.Range("W4").Value = 8.5
$ws.Range("X4").Value = 15
$ws.Range("Y4").Value = 12
$ws.Range("Z4").Value = 34
$ws.Range("AA4").Value = 29
$ws.Range("AB4").Value = 41
$ws.Range("AC4").Value = 8
$ws.Range("AD4").Value = 6
$ws.Range("AE4").Value = 17
$ws.Range("AF4").Value = 51
$ws.Range("AG4").Value = 7
$ws.Range("AH4").Value = 10
$ws.Range("AI4").Value = 9.5
$ws.Range("AJ4").Value = 21
$ws.Range("AK4").Value = 21
$ws.Range("AL4").Value = 34
$ws.Range("AM4").Value = 351
$ws.Range("AN4").Value = 5
$ws.Range("AO4").Value = 19
$ws.Range("AP4").Value = 29
$ws.Range("AQ4").Value = 67
$ws.Range("AR4").Value = 101
$ws.Range("AS4").Value = 251
$ws.Range("AT4").Value = 2.5
$ws.Range("AU4").Value = 8.5
$ws.Range("AV4").Value = 67
$ws.Range("AW4").Value = 4.33
$ws.Range("AX4").Value = 13
$ws.Range("AY4").Value = 26
$ws.Range("AZ4").Value = 41
$ws.Range("BA4").Value = 67
$ws.Range("BB4").Value = 201
$ws.Range("BC4").Value = 51
$ws.Range("BD4").Value = 51
